# Daily attendance processing - 2025-12-21 21:27:29
#
# The "Recorded By" column (G) lists the recorder(s) for each attendance
# session as a comma-separated string. For the specific recorder
# combinations below, rotate the list one position to the left (the
# first name moves to the end) to reflect the latest processing order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $text = $ws.Cells.Item($r, 7).Text

    if ($text -eq "system, backup@backdoor.com, System" -or `
        $text -eq "dnasr281@gmail.com, System" -or `
        $text -eq "System, dnasr281@gmail.com") {

        $parts = $text -split ", "
        $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
        $ws.Cells.Item($r, 7).Value = $rotated
    }
}
